{"js": "// Turn:\n//   \" ha az egyiken \u00e1ll). Lehet\u0151s\u00e9ge van eladni is mez\u0151it (\"\n// into:\n//   \" ha az egyiken \u00e1ll, max 5 h\u00e1z). Lehet\u0151s\u00e9ge van eladni is mez\u0151it (\"\n// by inserting \", max 5 h\u00e1z\" right after \"ha az egyiken \u00e1ll\" (and before\n// the closing paren that follows it).\nconst body = context.document.body;\nconst results = body.search(\"ha az egyiken \u00e1ll\", {\n  matchCase: false,\n  matchWholeWord: false\n});\nresults.load(\"text\");\nawait context.sync();\n\nif (results.items.length === 0) {\n  throw new Error(\"Target phrase 'ha az egyiken \u00e1ll' not found\");\n}\n\n// Insert the new text immediately at the end of the matched phrase, i.e.\n// right before the \")\" that follows it in the source paragraph.\nconst target = results.items[0];\ntarget.insertText(\", max 5 h\u00e1z\", Word.InsertLocation.end);\nawait context.sync();\n", "ps1": "# Turn:\n#   \" ha az egyiken \u00e1ll). Lehet\u0151s\u00e9ge van eladni is mez\u0151it (\"\n# into:\n#   \" ha az egyiken \u00e1ll, max 5 h\u00e1z). Lehet\u0151s\u00e9ge van eladni is mez\u0151it (\"\n# by replacing the phrase \"ha az egyiken \u00e1ll)\" with\n# \"ha az egyiken \u00e1ll, max 5 h\u00e1z)\" (keeping the closing paren as an anchor\n# so the match is unique and unambiguous).\n$d = $word.ActiveDocument\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \"ha az egyiken \u00e1ll)\"\n$find.Replacement.Text = \"ha az egyiken \u00e1ll, max 5 h\u00e1z)\"\n$find.Forward = $true\n$find.Wrap = 1  # wdFindContinue\n$find.Execute([ref]$null, [ref]$null, [ref]$null, [ref]$null, [ref]$null, [ref]$null, [ref]$null, [ref]$null, [ref]$null, [ref]$null, 2)  # wdReplaceAll\n"}
